# Added independent conspecifics corrections for allopatry and sympatry
#
# The slide's drawable content lives inside a single top-level group shape
# (Shapes.Item(1)); the five "rc23".."rc27" rectangles that need resizing
# are items 21-25 of that group's GroupItems collection.
#
# Left/Width are expressed in points (the PowerPoint COM unit); the literal
# values below were chosen so that, after the host's point->EMU conversion,
# the stored EMU values land exactly on the target offsets/extents.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)
$items = $g.GroupItems

# rc23 (id=23): width grows, position unchanged
$rc23 = $items.Item(21)
$rc23.Width = 124.61874389648438

# rc24 (id=24): width grows, position unchanged
$rc24 = $items.Item(22)
$rc24.Width = 22.241182327270508

# rc25 (id=25): left shifts right, width shrinks
$rc25 = $items.Item(23)
$rc25.Left = 296.93719482421875
$rc25.Width = 15.450000762939453

# rc26 (id=26): left shifts left, width grows substantially
$rc26 = $items.Item(24)
$rc26.Left = 242.26788330078125
$rc26.Width = 70.11921691894531

# rc27 (id=27): left shifts right slightly, width shrinks slightly
$rc27 = $items.Item(25)
$rc27.Left = 274.6959228515625
$rc27.Width = 37.69118118286133
